$wb = $excel.ActiveWorkbook
$wsFBS = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

# Update the refresh Timestamp (shared across all FBS rows, column AK)
for ($r = 2; $r -le 49; $r++) {
    $wsFBS.Range("AK$r").Value = "2024-11-04T16:21:42.443405"
}

$wsFBS.Range("N2").Value = "N"
$wsFBS.Range("O2").Value = 62.69
$wsFBS.Range("P2").Value = 20
$wsFBS.Range("Q2").Value = "N"
$wsFBS.Range("R2").Value = 3.95
$wsFBS.Range("U2").Value = 10.2
$wsFBS.Range("Y2").Value = 48.5
$wsFBS.Range("Z2").Value = -110
$wsFBS.Range("AE2").Value = -0.0202020202020202
$wsFBS.Range("M3").Value = "N"
$wsFBS.Range("O3").Value = 65.24000000000001
$wsFBS.Range("P3").Value = 17.1
$wsFBS.Range("R3").Value = 0.4
$wsFBS.Range("U3").Value = 5.7
$wsFBS.Range("Y3").Value = 48.5
$wsFBS.Range("Z3").Value = -110
$wsFBS.Range("AE3").Value = -0.0396039603960396
$wsFBS.Range("M4").Value = "WNW"
$wsFBS.Range("N4").Value = "WNW"
$wsFBS.Range("Q4").Value = "WNW"
$wsFBS.Range("M7").Value = "WSW"
$wsFBS.Range("N7").Value = "WSW"
$wsFBS.Range("O7").Value = 79.82000000000001
$wsFBS.Range("P7").Value = 14.1
$wsFBS.Range("Q7").Value = "WSW"
$wsFBS.Range("U7").Value = -0.3
$wsFBS.Range("Z7").Value = -115
$wsFBS.Range("M8").Value = "WNW"
$wsFBS.Range("N8").Value = "WNW"
$wsFBS.Range("Q8").Value = "WNW"
$wsFBS.Range("M10").Value = "W"
$wsFBS.Range("N10").Value = "W"
$wsFBS.Range("Q10").Value = "W"
$wsFBS.Range("Y10").Value = 57.5
$wsFBS.Range("AB10").Value = 3
$wsFBS.Range("AE10").Value = 0.03603603603603604
$wsFBS.Range("AF10").Value = 0
$wsFBS.Range("M11").Value = "ESE"
$wsFBS.Range("N11").Value = "ESE"
$wsFBS.Range("Q11").Value = "ESE"
$wsFBS.Range("AB11").Value = -26.5
$wsFBS.Range("AF11").Value = 0.5
$wsFBS.Range("M12").Value = "ENE"
$wsFBS.Range("N12").Value = "ENE"
$wsFBS.Range("Q12").Value = "ENE"
$wsFBS.Range("N15").Value = "ENE"
$wsFBS.Range("M16").Value = "W"
$wsFBS.Range("N16").Value = "WNW"
$wsFBS.Range("Q16").Value = "W"
$wsFBS.Range("AB16").Value = -1
$wsFBS.Range("AF16").Value = -0.5
$wsFBS.Range("N17").Value = "ENE"
$wsFBS.Range("Q17").Value = "ENE"
$wsFBS.Range("Z19").Value = -106
$wsFBS.Range("M20").Value = "ESE"
$wsFBS.Range("N20").Value = "ESE"
$wsFBS.Range("Q20").Value = "ESE"
$wsFBS.Range("Y20").Value = 48.5
$wsFBS.Range("Z20").Value = -110
$wsFBS.Range("AE20").Value = 0
$wsFBS.Range("M21").Value = "W"
$wsFBS.Range("N21").Value = "W"
$wsFBS.Range("Q21").Value = "W"
$wsFBS.Range("Y21").Value = 56.5
$wsFBS.Range("AB21").Value = -5
$wsFBS.Range("AE21").Value = 0
$wsFBS.Range("AF21").Value = 0.5
$wsFBS.Range("M22").Value = "WNW"
$wsFBS.Range("N22").Value = "WNW"
$wsFBS.Range("Q22").Value = "WNW"
$wsFBS.Range("AB22").Value = 11
$wsFBS.Range("AF22").Value = 1.5
$wsFBS.Range("M24").Value = "W"
$wsFBS.Range("N24").Value = "W"
$wsFBS.Range("Q24").Value = "W"
$wsFBS.Range("AB28").Value = 2.5
$wsFBS.Range("AF28").Value = -0.5
$wsFBS.Range("N29").Value = "WSW"
$wsFBS.Range("Q29").Value = "SW"
$wsFBS.Range("AB29").Value = -3
$wsFBS.Range("AF29").Value = -1
$wsFBS.Range("M30").Value = "WNW"
$wsFBS.Range("N30").Value = "WNW"
$wsFBS.Range("Q30").Value = "WNW"
$wsFBS.Range("M31").Value = "WNW"
$wsFBS.Range("N31").Value = "WNW"
$wsFBS.Range("Q31").Value = "WNW"
$wsFBS.Range("AB31").Value = 6.5
$wsFBS.Range("AF31").Value = 0.5
$wsFBS.Range("M33").Value = "ESE"
$wsFBS.Range("N33").Value = "ESE"
$wsFBS.Range("Q33").Value = "ESE"
$wsFBS.Range("AB35").Value = 10
$wsFBS.Range("AF35").Value = 0.5
$wsFBS.Range("M36").Value = "ENE"
$wsFBS.Range("N36").Value = "ENE"
$wsFBS.Range("Q36").Value = "ENE"
$wsFBS.Range("Y36").Value = 60.5
$wsFBS.Range("Z36").Value = -106
$wsFBS.Range("AB36").Value = -16.5
$wsFBS.Range("AE36").Value = 0.03418803418803419
$wsFBS.Range("AF36").Value = 3
$wsFBS.Range("M38").Value = "ESE"
$wsFBS.Range("N38").Value = "ESE"
$wsFBS.Range("Q38").Value = "ESE"
$wsFBS.Range("Y38").Value = 55.5
$wsFBS.Range("AB38").Value = -3
$wsFBS.Range("AE38").Value = -0.01769911504424779
$wsFBS.Range("AF38").Value = -2
$wsFBS.Range("N40").Value = "WSW"
$wsFBS.Range("Q40").Value = "SW"
$wsFBS.Range("AB40").Value = -25
$wsFBS.Range("AF40").Value = 0.5
$wsFBS.Range("M41").Value = "WSW"
$wsFBS.Range("N41").Value = "WSW"
$wsFBS.Range("Q41").Value = "WSW"
$wsFBS.Range("AB41").Value = -24
$wsFBS.Range("AF41").Value = -0.5
$wsFBS.Range("M42").Value = "ENE"
$wsFBS.Range("AB42").Value = 3
$wsFBS.Range("AF42").Value = -0.5
$wsFBS.Range("AB43").Value = -2
$wsFBS.Range("AF43").Value = -0.5
$wsFBS.Range("M44").Value = "WSW"
$wsFBS.Range("N44").Value = "WSW"
$wsFBS.Range("Q44").Value = "WSW"
$wsFBS.Range("AB44").Value = -7.5
$wsFBS.Range("AF44").Value = 0.5
$wsFBS.Range("M45").Value = "ESE"
$wsFBS.Range("N45").Value = "ESE"
$wsFBS.Range("Q45").Value = "ESE"
$wsFBS.Range("AB45").Value = -24.5
$wsFBS.Range("AF45").Value = -1
$wsFBS.Range("M46").Value = "W"
$wsFBS.Range("N46").Value = "W"
$wsFBS.Range("Q46").Value = "W"
$wsFBS.Range("M47").Value = "ESE"
$wsFBS.Range("N47").Value = "ESE"
$wsFBS.Range("Q47").Value = "ESE"
$wsFBS.Range("M48").Value = "W"
$wsFBS.Range("N48").Value = "WNW"
$wsFBS.Range("Q48").Value = "WNW"
$wsOther.Range("O2").Value = "WNW"
$wsOther.Range("P2").Value = "WNW"
$wsOther.Range("S2").Value = "WNW"
$wsOther.Range("O3").Value = "WNW"
$wsOther.Range("P3").Value = "WNW"
$wsOther.Range("S3").Value = "WNW"
$wsOther.Range("O6").Value = "W"
$wsOther.Range("P6").Value = "W"
$wsOther.Range("S6").Value = "W"
$wsOther.Range("O9").Value = "ESE"
$wsOther.Range("P9").Value = "ESE"
$wsOther.Range("S9").Value = "ESE"
$wsOther.Range("O10").Value = "WSW"
$wsOther.Range("P10").Value = "WSW"
$wsOther.Range("S10").Value = "WSW"
$wsOther.Range("O12").Value = "W"
$wsOther.Range("P12").Value = "W"
$wsOther.Range("S12").Value = "W"
$wsOther.Range("O16").Value = "W"
$wsOther.Range("P16").Value = "W"
$wsOther.Range("S16").Value = "W"
$wsOther.Range("O18").Value = "WNW"
$wsOther.Range("P18").Value = "W"
$wsOther.Range("S18").Value = "W"
$wsOther.Range("O21").Value = "WSW"
$wsOther.Range("P21").Value = "WSW"
$wsOther.Range("S21").Value = "WSW"
$wsOther.Range("O30").Value = "WSW"
$wsOther.Range("O31").Value = "W"
$wsOther.Range("P31").Value = "W"
$wsOther.Range("S31").Value = "W"
$wsOther.Range("P35").Value = "WNW"
$wsOther.Range("S35").Value = "WNW"
$wsOther.Range("O37").Value = "ENE"
$wsOther.Range("P37").Value = "ENE"
$wsOther.Range("S37").Value = "ENE"
$wsOther.Range("O38").Value = "ENE"
$wsOther.Range("P38").Value = "ENE"
$wsOther.Range("S38").Value = "ENE"
$wsOther.Range("O39").Value = "ESE"
$wsOther.Range("S40").Value = "SW"
$wsOther.Range("O42").Value = "W"
$wsOther.Range("P42").Value = "W"
$wsOther.Range("S42").Value = "W"
$wsOther.Range("P43").Value = "W"
$wsOther.Range("S43").Value = "WNW"
$wsOther.Range("S45").Value = "SE"
$wsOther.Range("O46").Value = "ENE"
$wsOther.Range("P46").Value = "ENE"
$wsOther.Range("S46").Value = "ENE"
